# Updated CSD util and testbase
# - Populate the "Incident" sheet (2nd sheet) with IncName/contractName test rows
# - Match header styling used on the "Events" sheet (bold + yellow fill)
# - Size the columns to fit the new content
# - Set page setup (paper size / orientation) to match the rest of the workbook
# - Make "Incident" the active/selected sheet & select cell A3 on it

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # "Incident" sheet

# --- Data -------------------------------------------------------------
$ws2.Range("A1").Value = "IncName"
$ws2.Range("B1").Value = "contractName"
$ws2.Range("A2").Value = "LDAP is Down.This is a Test Incident"
$ws2.Range("B2").Value = "CAH_TEST"
$ws2.Range("A3").Value = "Sentinal is Down.This is a Test Incident"
$ws2.Range("B3").Value = "ES_TEST"

# --- Header formatting (bold text on yellow fill, like sheet "Events") -
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("A1:B1").Interior.Color = 65535

# --- Column sizing ------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 32.333333333333336
$ws2.Columns.Item(2).ColumnWidth = 11.666666666666666

# --- Page setup (match paper size / orientation used elsewhere) --------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Selection / active sheet -------------------------------------------
$ws2.Range("A3").Select() | Out-Null
$ws2.Activate()
